$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1, matching the same formatting as the other headers
# (bold font, centered/top alignment, thin box border)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data values for I2:I58 and J2:J58
$iValues = @(6,8,7,7,5,10,9,6,8,7,4,8,7,9,6,7,4,6,9,7,6,6,6,6,1,6,10,8,9,7,7,6,10,7,8,7,6,7,6,7,8,6,8,8,8,5,9,8,6,3,8,6,6,5,4,7,5)
$jValues = @(6,8,7,7,6,10,9,6,8,7,4,8,8,9,7,7,4,6,9,8,6,6,7,7,1,6,10,8,9,7,7,6,11,7,8,8,7,7,6,7,8,6,8,8,8,6,9,8,6,3,8,6,6,5,4,7,5)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
